$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.267.75"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.876.02"
$ws.Range("E3").Value = "  +3.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.35"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5026"
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3949"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09904"
$ws.Range("E9").Value = "  +27.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.144"
$ws.Range("E10").Value = "  +3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.29"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.491"
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.04"
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("D14").Value = "1.867.49"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.000"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.407"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("E17").Value = "  +5.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.69"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06646"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.46"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.135"
$ws.Range("E22").Value = "  +2.37%  "
$ws.Range("D23").Value = "28.320.69"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.36"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.574"
$ws.Range("E26").Value = "  +4.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.33"
$ws.Range("E27").Value = "  +4.37%  "
$ws.Range("D28").Value = "2.086.34"
$ws.Range("E28").Value = "  +3.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.28"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.88"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1064"
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.065"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.648"
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.607"
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06832"
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.541"
$ws.Range("E36").Value = "  +4.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02403"
$ws.Range("E37").Value = "  +2.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2188"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.54"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.031"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6326"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.176"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.45"
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6033"
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.668"
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.271"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.12"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.997"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.203"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.127"
$ws.Range("E51").Value = "  +6.13%  "
